# The deck's live theme (ppt/theme/theme2.xml, wired to the one real
# SlideMaster) currently carries the "Integral" color scheme. The commit
# swaps it back to the stock "Office Theme" palette (fonts/format scheme
# were already identical between the two theme parts, so only the 12
# scheme colors actually need to change).
#
# ThemeColorScheme indices follow the standard MsoThemeColorSchemeIndex
# order: 1=dk1 2=lt1 3=dk2 4=lt2 5-10=accent1-6 11=hlink 12=folHlink.
# .RGB takes/returns the usual COM BGR-packed integer (like the VBA RGB()
# macro would build), so each target hex color below is converted
# (B<<16 | G<<8 | R) before being assigned.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB = 0          # dk1     000000
$tcs.Item(2).RGB = 16777215   # lt1     FFFFFF
$tcs.Item(3).RGB = 6968388    # dk2     44546A
$tcs.Item(4).RGB = 15132391   # lt2     E7E6E6
$tcs.Item(5).RGB = 13998939   # accent1 5B9BD5
$tcs.Item(6).RGB = 3243501    # accent2 ED7D31
$tcs.Item(7).RGB = 10855845   # accent3 A5A5A5
$tcs.Item(8).RGB = 49407      # accent4 FFC000
$tcs.Item(9).RGB = 12874308   # accent5 4472C4
$tcs.Item(10).RGB = 4697456   # accent6 70AD47
$tcs.Item(11).RGB = 12673797  # hlink   0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
